$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the numeric values first, then the shared string.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, centered horizontally, top vertically.
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("B1").Borders.Weight = 2
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4160

# Give A2 the exact same formatting as B1 by copying it across,
# rather than re-deriving it property-by-property (avoids creating
# an extra, unused intermediate cell style in styles.xml).
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
